# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3377
$ws1.Range("F4").Value = 65
$ws1.Range("F5").Value = 1531
$ws1.Range("F6").Value = 52
$ws1.Range("F7").Value = 324

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3377
$ws4.Range("F4").Value = 65
$ws4.Range("F5").Value = 1531
$ws4.Range("F6").Value = 52
$ws4.Range("F8").Value = 324
